$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2. Информация об организации -- update contact details for the indicator's
# responsible organization (department renamed + new contact person/email/phone/site).
# Written in this order so the shared-string table append order matches the
# authored workbook (phone, site, email, contact person, org name).
$ws.Range("B9").Value  = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B8").Value  = "yryskan.kalymbetova@gmail.com "
$ws.Range("B7").Value  = "Калымбетова Ы.И."
$ws.Range("B6").Value  = "Национальный статистический комитет КР `n(Управление статистики домашних хозяйств)"

# Move the active selection to B4, matching the authored workbook's saved view state.
$ws.Range("B4").Select() | Out-Null
